# Refresh cryptos list values (Price / Volume(1h), and the two swapped rows)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.125.47"
$ws.Range("E2").Value = "  -3.70%  "
$ws.Range("D3").Value = "1.849.14"
$ws.Range("E3").Value = "  -2.70%  "
$ws.Range("D4").Value = "'0.9993"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "'0.7071"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -5.52%  "
$ws.Range("D6").Value = "'238.25"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -1.91%  "
$ws.Range("D7").Value = "'0.9994"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +0.18%  "
$ws.Range("D8").Value = "'0.3053"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -4.00%  "
$ws.Range("D9").Value = "'0.07506"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +3.27%  "
$ws.Range("D10").Value = "'23.40"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -7.06%  "
$ws.Range("D11").Value = "'0.08137"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -2.87%  "
$ws.Range("D12").Value = "1.870.61"
$ws.Range("E12").Value = "  -1.22%  "
$ws.Range("E13").Value = "  -5.30%  "
$ws.Range("D14").Value = "'5.225"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -4.46%  "
$ws.Range("D15").Value = "'89.23"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -4.34%  "
$ws.Range("D16").Value = "29.149.64"
$ws.Range("E16").Value = "  -3.46%  "
$ws.Range("D17").Value = "'5.793"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -6.49%  "
$ws.Range("D18").Value = "'239.85"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -4.74%  "
$ws.Range("E19").Value = "  -2.76%  "
$ws.Range("D20").Value = "'13.08"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -4.68%  "
$ws.Range("D21").Value = "'1.000"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +0.31%  "
$ws.Range("D22").Value = "2.096.71"
$ws.Range("E22").Value = "  -2.15%  "
$ws.Range("D23").Value = "'0.9994"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +0.05%  "
$ws.Range("D24").Value = "'7.561"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -6.03%  "
$ws.Range("D25").Value = "'0.1467"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -8.32%  "
$ws.Range("D26").Value = "'8.987"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -3.83%  "
$ws.Range("D27").Value = "'161.30"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -2.10%  "
$ws.Range("E28").Value = "  -4.46%  "
$ws.Range("D29").Value = "'1.943"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -6.46%  "
$ws.Range("D30").Value = "'1.386"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -6.26%  "
$ws.Range("D31").Value = "'4.548"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -1.79%  "
$ws.Range("D32").Value = "'1.495"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -2.84%  "
$ws.Range("D33").Value = "'4.004"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -5.63%  "
$ws.Range("D34").Value = "'0.05175"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -4.54%  "
$ws.Range("D35").Value = "'1.188"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -6.04%  "
$ws.Range("D36").Value = "'1.036"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +4.03%  "
$ws.Range("D37").Value = "'0.7061"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -7.86%  "
$ws.Range("D38").Value = "'2.642"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -2.59%  "
$ws.Range("E39").Value = "  -5.81%  "
$ws.Range("D40").Value = "'2.678"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -3.45%  "
$ws.Range("D41").Value = "'0.9339"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +6.85%  "
$ws.Range("D42").Value = "'5.997"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -1.63%  "
$ws.Range("D43").Value = "1.078.28"
$ws.Range("E43").Value = "  -1.96%  "
$ws.Range("D44").Value = "'0.4307"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -6.42%  "
$ws.Range("D45").Value = "'70.20"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -3.86%  "
$ws.Range("E46").Value = "  -0.10%  "
$ws.Range("D47").Value = "'102.27"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -2.22%  "
$ws.Range("E48").Value = "  -6.78%  "
$ws.Range("B49").Value = "RocketPoolETH"
$ws.Range("C49").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D49").Value = "1.993.15"
$ws.Range("E49").Value = "  -2.66%  "
$ws.Range("B50").Value = "Aptos"
$ws.Range("C50").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D50").Value = "'7.070"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -7.62%  "
$ws.Range("E51").Value = "  -4.86%  "
